$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data layout: Row, D(Fecha serial), J(Volumen), K(Precio minimo), L(Precio maximo), M(Precio promedio ponderado), P(Precio $/Kg)
$data = @(
    @(2,44330,30,15000,15000,15000,1000),
    @(3,44315,65,14000,15000,14538,969),
    @(4,44729,85,16000,17000,16529,1102),
    @(5,44316,45,14000,15000,14444,963),
    @(6,44326,45,15000,15000,15000,1000),
    @(7,44756,50,15000,15000,15000,1000),
    @(8,44319,50,15000,15000,15000,1000),
    @(9,44406,50,22000,22000,22000,1467),
    @(10,44750,85,15000,16000,15471,1031),
    @(11,44760,105,15000,16000,15524,1035),
    @(12,44322,70,14000,15000,14500,967),
    @(13,44715,85,15000,15500,15235,1016),
    @(14,44727,60,15000,15000,15000,1000),
    @(15,44309,50,15000,15000,15000,1000),
    @(16,44742,85,15000,16000,15529,1035),
    @(17,44719,60,15000,15000,15000,1000),
    @(18,44720,85,15000,16000,15529,1035),
    @(19,44755,100,15000,16000,15550,1037),
    @(20,44722,95,15000,15500,15263,1018),
    @(21,44320,40,15000,15000,15000,1000),
    @(22,44448,85,21000,22000,21529,1435),
    @(23,44411,50,22000,22000,22000,1467),
    @(24,44714,100,15000,15500,15250,1017),
    @(25,44323,40,15000,15000,15000,1000),
    @(26,44308,40,16000,16000,16000,1067),
    @(27,44329,35,15000,15000,15000,1000),
    @(28,44312,80,13000,14000,13562,904),
    @(29,44455,35,22000,22000,22000,1467),
    @(30,44334,50,14000,14000,14000,933),
    @(31,44343,40,15000,15000,15000,1000),
    @(32,44721,130,14000,15000,14500,967),
    @(33,44333,35,15000,15000,15000,1000),
    @(34,44753,80,15000,16000,15500,1033),
    @(35,44749,100,16000,17000,16450,1097),
    @(36,44344,40,20000,20000,20000,1333),
    @(37,44754,50,15000,15000,15000,1000),
    @(38,44725,85,14000,15000,14471,965),
    @(39,44736,82,16000,17000,16488,1099),
    @(40,44314,45,15000,15000,15000,1000),
    @(41,44313,40,14000,14000,14000,933),
    @(42,44327,35,15000,15000,15000,1000),
    @(43,44746,103,15000,16000,15563,1038),
    @(44,44328,38,15000,15000,15000,1000),
    @(45,44452,73,22000,23000,22479,1499),
    @(46,44341,40,15000,15000,15000,1000),
    @(47,44370,50,18000,18000,18000,1200),
    @(48,44748,73,15000,16000,15521,1035),
    @(49,44726,55,15000,15000,15000,1000),
    @(50,44336,65,14000,15000,14462,964),
    @(51,44340,47,14000,14000,14000,933),
    @(52,44399,38,22000,22000,22000,1467),
    @(53,44377,80,18000,19000,18500,1233),
    @(54,44397,73,21000,22000,21521,1435),
    @(55,44747,40,16000,16000,16000,1067),
    @(56,44757,40,15000,15000,15000,1000),
    @(57,44321,38,15000,15000,15000,1000),
    @(58,44438,75,19000,20000,19467,1298)
)

foreach ($row in $data) {
    $r = $row[0]
    $d = $row[1]
    $j = $row[2]
    $k = $row[3]
    $l = $row[4]
    $m = $row[5]
    $p = $row[6]

    $cellD = $ws.Cells.Item($r, 4)
    $cellD.Value = $d
    $cellD.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 10).Value = $j
    $ws.Cells.Item($r, 11).Value = $k
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 16).Value = $p
}

# New row 58 needs the full set of constant columns (identical across every data row in this sheet)
$ws.Cells.Item(58, 1).Value = 3
$ws.Cells.Item(58, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(58, 3).Value = "Coquimbo"
$ws.Cells.Item(58, 5).Value = 5
$ws.Cells.Item(58, 6).Value = 100112035
$ws.Cells.Item(58, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(58, 8).Value = "Sin especificar"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(58, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(58, 17).Value = 15
$ws.Cells.Item(58, 18).Value = "Hortaliza"

